# Add a "Save" column (H) to the s_vals worksheet.
# H1 gets the header "Save" (styled like the other header cells).
# For each data row, H = 1 when the "sum" value in column G is a save-type
# outing (G > 20), otherwise H = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching the style used by the other header cells (B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$lastRow = 64
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value()
    if ($g -gt 20) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
